$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Update Sheet8: add a new value in H1, and change the selection to H3.
# ---------------------------------------------------------------------------
$ws8 = $wb.Worksheets.Item("Sheet8")
$ws8.Range("H1").Value = 0
$ws8.Range("H3").Select()

# ---------------------------------------------------------------------------
# 2. Add a new worksheet ("Sheet9") after the last existing sheet.
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws9 = $wb.Worksheets.Add($null, $lastSheet)

# Row 2
$ws9.Range("B2").Value = 28
$ws9.Range("C2").Value = 63
$ws9.Range("D2").Value = 63
$ws9.Range("E2").Value = 63
$ws9.Range("F2").Value = 91
$ws9.Range("G2").Value = 91

# Row 4
$ws9.Range("B4").Value = 0
$ws9.Range("C4").Value = 28
$ws9.Range("D4").Value = 56
$ws9.Range("E4").Value = 87
$ws9.Range("F4").Value = 56

# Row 5
$ws9.Range("B5").Value = 32
$ws9.Range("C5").Value = 70
$ws9.Range("D5").Value = 70
$ws9.Range("E5").Value = 70
$ws9.Range("F5").Value = 70
$ws9.Range("G5").Value = 99
$ws9.Range("H5").Value = 70
$ws9.Range("I5").Value = 70

# Row 7
$ws9.Range("C7").Value = 70
$ws9.Range("E7").Value = 0

# Row 9
$ws9.Range("E9").Value = 0

# Row 10
$ws9.Range("A10").Value = 34
$ws9.Range("C10").Value = 29

# Selection + active sheet for the new sheet.
$ws9.Range("C10").Select()
$ws9.Activate()

# ---------------------------------------------------------------------------
# 3. Window chrome: minimize the workbook window and scroll the tab strip so
#    that the first visible tab is Sheet5 (index 2).
# ---------------------------------------------------------------------------
$win = $wb.Windows.Item(1)
try { $win.WindowState = -4140 } catch { }
try { $win.ScrollWorkbookTabs(2) } catch { }
